$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row permutation: for each destination row, the source row whose
# D, L, M, N, O, P, S values are copied into it.
$mapping = @{
    2  = 3
    3  = 6
    4  = 11
    5  = 12
    6  = 2
    7  = 8
    8  = 9
    9  = 4
    10 = 5
    11 = 7
    12 = 10
}

# Snapshot the original values of the affected columns before overwriting
# anything, so that later writes don't clobber values still needed as a
# source for another row.
$orig = @{}
foreach ($r in 2..12) {
    $orig[$r] = @{
        D = $ws.Range("D$r").Value2
        L = $ws.Range("L$r").Value2
        M = $ws.Range("M$r").Value2
        N = $ws.Range("N$r").Value2
        O = $ws.Range("O$r").Value2
        P = $ws.Range("P$r").Value2
        S = $ws.Range("S$r").Value2
    }
}

foreach ($r in 2..12) {
    $src = $mapping[$r]
    $vals = $orig[$src]
    $ws.Range("D$r").Value2 = $vals.D
    $ws.Range("L$r").Value2 = $vals.L
    $ws.Range("M$r").Value2 = $vals.M
    $ws.Range("N$r").Value2 = $vals.N
    $ws.Range("O$r").Value2 = $vals.O
    $ws.Range("P$r").Value2 = $vals.P
    $ws.Range("S$r").Value2 = $vals.S
}
